# Auto-generated Excel COM-interop edit script
# Applies the diff: updates Timestamp + several odds/weather cell values
$wb = $excel.ActiveWorkbook
$wsFBS = $wb.Worksheets.Item("FBS")
$wsOther = $wb.Worksheets.Item("Other")

# --- Update Timestamp shared string (AK column, FBS sheet, rows 2-50) ---
# Every row previously showed 2024-10-09T05:15:54.232776; update every one
# so they all share the single updated timestamp value, matching the commit message.
for ($row = 2; $row -le 50; $row++) {
    $wsFBS.Cells.Item($row, 37).Value = "2024-10-09T16:21:43.374041"  # column AK = 37
}

# --- FBS sheet: numeric cell updates ---
$wsFBS.Range("AB2").Value = 2
$wsFBS.Range("AF2").Value = -2
$wsFBS.Range("Z5").Value = -105
$wsFBS.Range("Y6").Value = 47.5
$wsFBS.Range("Z6").Value = -110
$wsFBS.Range("AB6").Value = 3.5
$wsFBS.Range("AE6").Value = -0.04040404040404041
$wsFBS.Range("AF6").Value = 0.5
$wsFBS.Range("Y7").Value = 50.5
$wsFBS.Range("Z7").Value = -112
$wsFBS.Range("AE7").Value = -0.01941747572815534
$wsFBS.Range("Z10").Value = -105
$wsFBS.Range("AB10").Value = 7
$wsFBS.Range("AF10").Value = -0.5
$wsFBS.Range("AB12").Value = -20
$wsFBS.Range("AF12").Value = 1.5
$wsFBS.Range("Z13").Value = -105
$wsFBS.Range("O14").Value = 81.74000000000001
$wsFBS.Range("P14").Value = 10
$wsFBS.Range("S14").Value = -0.22
$wsFBS.Range("T14").Value = -0.22
$wsFBS.Range("U14").Value = -3.5
$wsFBS.Range("Y14").Value = 60.5
$wsFBS.Range("Z14").Value = -106
$wsFBS.Range("AE14").Value = -0.01626016260162602
$wsFBS.Range("Z15").Value = -115
$wsFBS.Range("O16").Value = 69.70999999999999
$wsFBS.Range("P16").Value = 8.6
$wsFBS.Range("U16").Value = 4.3
$wsFBS.Range("Z16").Value = -108
$wsFBS.Range("AB16").Value = -22
$wsFBS.Range("AF16").Value = 1
$wsFBS.Range("Z18").Value = -122
$wsFBS.Range("AB18").Value = -4.5
$wsFBS.Range("AF18").Value = 0
$wsFBS.Range("AB20").Value = -11
$wsFBS.Range("AF20").Value = 1
$wsFBS.Range("Z21").Value = -108
$wsFBS.Range("Z22").Value = -105
$wsFBS.Range("Y24").Value = 43.5
$wsFBS.Range("AE24").Value = -0.06451612903225806
$wsFBS.Range("Y25").Value = 59.5
$wsFBS.Range("Z25").Value = -105
$wsFBS.Range("AE25").Value = 0.03478260869565217
$wsFBS.Range("Y27").Value = 50.5
$wsFBS.Range("Z27").Value = -105
$wsFBS.Range("AE27").Value = -0.01941747572815534
$wsFBS.Range("Y28").Value = 41.5
$wsFBS.Range("Z28").Value = -104
$wsFBS.Range("AE28").Value = 0.05063291139240506
$wsFBS.Range("Y29").Value = 46.5
$wsFBS.Range("Z29").Value = -110
$wsFBS.Range("AE29").Value = 0.02197802197802198
$wsFBS.Range("Y30").Value = 58.5
$wsFBS.Range("AE30").Value = -0.03305785123966942
$wsFBS.Range("Y32").Value = 52.5
$wsFBS.Range("Z32").Value = -115
$wsFBS.Range("AE32").Value = 0
$wsFBS.Range("Z34").Value = -102
$wsFBS.Range("Y35").Value = 58.5
$wsFBS.Range("Z35").Value = -108
$wsFBS.Range("AE35").Value = -0.03305785123966942
$wsFBS.Range("Z38").Value = -110
$wsFBS.Range("Y39").Value = 49.5
$wsFBS.Range("Z39").Value = -106
$wsFBS.Range("AE39").Value = -0.0198019801980198
$wsFBS.Range("Y41").Value = 56.5
$wsFBS.Range("Z41").Value = -115
$wsFBS.Range("AE41").Value = 0.03669724770642202
$wsFBS.Range("Y46").Value = 53.5
$wsFBS.Range("Z46").Value = -110
$wsFBS.Range("AE46").Value = 0.03883495145631068
$wsFBS.Range("Z48").Value = -122

# --- FBS sheet: wind-direction (string) cell updates ---
$wsFBS.Range("M16").Value = "SW"
$wsFBS.Range("Q20").Value = "ENE"
$wsFBS.Range("Q39").Value = "NNE"
$wsFBS.Range("Q43").Value = "NE"
$wsFBS.Range("Q49").Value = "SE"
$wsFBS.Range("Q50").Value = "WNW"

# --- Other sheet: wind-direction (string) cell updates ---
$wsOther.Range("S14").Value = "NNE"
$wsOther.Range("S36").Value = "NE"

Write-Host "Applied all cell updates"